## Add new team member "Sadie Drouin" to the people workbook.
##
## Sheet2 is a "most recently added" staging row that always mirrors the
## newest entry appended to the bottom of Sheet1. This edit:
##   1. Appends Sadie Drouin as a brand-new row at the bottom of Sheet1 (row 36)
##      in the "Data Generation" category, with a mailto hyperlink on her email.
##   2. Replaces the Sheet2 staging row (previously Haley Fritch) with Sadie's
##      info as well, using the same layout (first, last, email, long_bio,
##      category) -- no importance/title columns for her.
##   3. Makes Sheet2 the active/selected sheet (previously Sheet1 was active).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$firstName = "Sadie"
$lastName = "Drouin"
$email = "drouin@broadinstitute.org"
$longBio = "Sadie is a Research Associate I for the McCarroll & Macosko BICAN project. She is a recent graduate from Wheaton College MA with a B.S. in Neuroscience on the Pre-Medical Track. Prior to joining the BICAN project, in her Behavioral Neuroscience lab she completed a senior honors thesis titled: The effects of maternal separation and social isolation on memory and myelin in adolescent rats. In addition to being a research associate, she is also a medical assistant in primary care and sports medicine at DMC primary care. "
$category = "Data Generation"

# --- Sheet1: append Sadie as new last row (row 36) ---
$newRow = 36
$ws1.Cells.Item($newRow, 1).Value = $firstName
$ws1.Cells.Item($newRow, 2).Value = $lastName
$ws1.Cells.Item($newRow, 7).Value = $longBio
$ws1.Cells.Item($newRow, 6).Value = $email
$ws1.Cells.Item($newRow, 8).Value = $category

$ws1.Hyperlinks.Add($ws1.Cells.Item($newRow, 6), "mailto:" + $email, "", "", $email) | Out-Null
$ws1.Cells.Item($newRow, 6).Style = "Hyperlink"

$ws1.Range("A36:H36").Select() | Out-Null

# --- Sheet2: clear existing staging row and replace with Sadie's data ---
$ws2.Range("A2:H2").ClearContents() | Out-Null
$ws2.Cells.Item(2, 1).Value = $firstName
$ws2.Cells.Item(2, 2).Value = $lastName
$ws2.Cells.Item(2, 7).Value = $longBio
$ws2.Cells.Item(2, 6).Value = $email
$ws2.Cells.Item(2, 8).Value = $category

$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 6), "mailto:" + $email, "", "", $email) | Out-Null
$ws2.Cells.Item(2, 6).Style = "Hyperlink"

# --- Sheet2 becomes the active tab (select it last so it "wins") ---
$ws2.Select()
